# Applies the "additional scraping" commit:
#  1. Inserts a new "Player Info" sheet at the front with ID/NAME/BATTING_HAND/BOWL_STYLE.
#  2. Keeps "ODI Batting" and "ODI Bowling" (now 2nd/3rd tabs) but:
#       - renames MATCH_CARD_LINK -> MATCH_CODE
#       - replaces the full scorecard URL with just the numeric match code
#  3. Appends a new "ODI Batting Extra" sheet with per-match batting detail.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

function Style-Header($rng) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# 1. "Player Info" sheet - brand new, becomes the first tab.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le 4; $c++) {
    Set-TextCell $playerInfo 1 $c $piHeaders[$c - 1]
}
Style-Header ($playerInfo.Range("A1:D1"))

$piRow = @("4238", "Christopher James Jordan", "Right Handed", "Right Arm Fast Medium")
for ($c = 1; $c -le 4; $c++) {
    Set-TextCell $playerInfo 2 $c $piRow[$c - 1]
}

# ---------------------------------------------------------------------------
# 2. "ODI Batting" - rename MATCH_CARD_LINK column, strip URLs to bare codes.
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
Set-TextCell $batting 1 4 "MATCH_CODE"

$lastRow = $batting.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $link = $batting.Cells.Item($r, 4).Text
    $code = $link -replace '^.*MatchCode=', ''
    Set-TextCell $batting $r 4 $code

    # Rows where INNING_NUMBER (column B) is blank ("did not bat") used to
    # carry a stray empty inline-string cell; clear it so the cell is
    # dropped from the sheet entirely, matching the regenerated export.
    $inning = $batting.Cells.Item($r, 2).Text
    if ([string]::IsNullOrEmpty($inning)) {
        $batting.Cells.Item($r, 2).ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" - same MATCH_CARD_LINK -> MATCH_CODE treatment.
# ---------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
Set-TextCell $bowling 1 2 "MATCH_CODE"

$lastRowBowl = $bowling.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowBowl; $r++) {
    $link = $bowling.Cells.Item($r, 2).Text
    $code = $link -replace '^.*MatchCode=', ''
    Set-TextCell $bowling $r 2 $code
}

# ---------------------------------------------------------------------------
# 4. "ODI Batting Extra" - brand new, appended as the last tab.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le 6; $c++) {
    Set-TextCell $extra 1 $c $exHeaders[$c - 1]
}
Style-Header ($extra.Range("A1:F1"))

# MATCH_CODE, BATTING_POSITION (numeric or blank), NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$exRows = @(
    @("3712", $null, "", "", "", "NO"),
    @("3714", 9, "", "", "", "NO"),
    @("3716", 9, "1", "0", "1.98%", "NO"),
    @("3718", 9, "0", "0", "", "NO"),
    @("3727", 9, "1", "0", "7.26%", "NO"),
    @("3780", 9, "0", "0", "", "NO"),
    @("3785", 8, "", "", "", "YES"),
    @("3804", 9, "0", "0", "0.49%", "NO"),
    @("3805", 10, "0", "0", "0.82%", "NO"),
    @("3885", 8, "0", "1", "3.51%", "NO"),
    @("3887", 8, "", "", "", "NO"),
    @("3889", $null, "", "", "", "NO"),
    @("3908", 8, "", "", "", "NO"),
    @("3911", 9, "", "", "", "NO"),
    @("3930", $null, "", "", "", "NO"),
    @("3932", 10, "0", "1", "4.97%", "NO"),
    @("4401", $null, "", "", "", "NO"),
    @("4405", $null, "", "", "", "NO"),
    @("4408", $null, "", "", "", "NO"),
    @("4660", $null, "", "", "", "NO")
)

$r = 2
foreach ($row in $exRows) {
    Set-TextCell $extra $r 1 $row[0]

    $posCell = $extra.Cells.Item($r, 2)
    if ($null -eq $row[1]) {
        # BATTING_POSITION unknown for this match - keep the cell present
        # (matches the regenerated export's empty <c> placeholder) but
        # General-formatted like the numeric cells around it.
        $posCell.NumberFormat = "General"
        $posCell.Value = ""
    } else {
        $posCell.NumberFormat = "General"
        $posCell.Value = $row[1]
    }

    Set-TextCell $extra $r 3 $row[2]
    Set-TextCell $extra $r 4 $row[3]
    Set-TextCell $extra $r 5 $row[4]
    Set-TextCell $extra $r 6 $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 5. Tab order + active sheet: Player Info, ODI Batting, ODI Bowling,
#    ODI Batting Extra - with Player Info left active/selected.
# ---------------------------------------------------------------------------
$playerInfo.Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item(1).Activate()

Write-Host "edit complete"
